# "Swarm missions, data tweaks"
#
# 1. Switch the active sheet from "Item Levels" (index 3) to
#    "Levels and Experience" (index 1) -- moves workbookView/activeTab
#    and the per-sheet tabSelected flag, and updates each sheet's
#    last-used selection.
# 2. Tweak the level-scaling base multiplier on "Levels and Experience"
#    (G2: 1.65 -> 1.62); all of the dependent INT/ROUND formulas in
#    columns C and D recompute automatically.

$wb = $excel.ActiveWorkbook

$wsItemLevels = $wb.Worksheets.Item("Item Levels")
$wsLevelsExp  = $wb.Worksheets.Item("Levels and Experience")

# Record the previous selection on "Item Levels" before it loses focus,
# then leave a fresh selection there (matches the post-edit snapshot).
$wsItemLevels.Range("J7").Select()

# Make "Levels and Experience" the active sheet (activeTab 3 -> 1) and
# set its new selection.
$wsLevelsExp.Activate()
$wsLevelsExp.Range("G2").Value = 1.62
$wsLevelsExp.Range("G3").Select()
